# Applies the cryptos.xlsx data refresh described by the commit diff.
# All written values are plain text (prices/URLs/percent strings), so we
# force text storage via NumberFormat "@" before writing, then clear the
# format again so the cell keeps its original (unstyled) appearance -
# otherwise Excel auto-coerces numeric-looking strings like "0.999" into
# real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "65.448.76"
Set-TextValue "E2" "  -1.83%  "
Set-TextValue "D3" "3.515.56"
Set-TextValue "E3" "  -2.29%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.13%  "
Set-TextValue "D5" "601.90"
Set-TextValue "E5" "  -1.20%  "
Set-TextValue "D6" "143.74"
Set-TextValue "E6" "  -1.92%  "
Set-TextValue "D7" "3.514.09"
Set-TextValue "E7" "  -2.39%  "
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.15%  "
Set-TextValue "D9" "0.513"
Set-TextValue "E9" "  +4.54%  "
Set-TextValue "E10" "  -2.40%  "
Set-TextValue "D11" "7.84"
Set-TextValue "E11" "  -1.40%  "
Set-TextValue "D12" "0.404"
Set-TextValue "E12" "  -2.82%  "
Set-TextValue "D13" "4.108.78"
Set-TextValue "E13" "  -2.45%  "
Set-TextValue "E14" "  -5.17%  "
Set-TextValue "D15" "28.50"
Set-TextValue "E15" "  -5.21%  "
Set-TextValue "D16" "3.514.44"
Set-TextValue "E16" "  -3.40%  "
Set-TextValue "E17" "  +1.47%  "
Set-TextValue "D18" "65.395.94"
Set-TextValue "E18" "  -2.08%  "
Set-TextValue "D19" "11.06"
Set-TextValue "E19" "  -3.74%  "
Set-TextValue "E20" "  -1.27%  "
Set-TextValue "D21" "14.33"
Set-TextValue "E21" "  -4.81%  "
Set-TextValue "D22" "416.84"
Set-TextValue "E22" "  -3.76%  "
Set-TextValue "D23" "0.598"
Set-TextValue "E23" "  -4.32%  "
Set-TextValue "D24" "77.56"
Set-TextValue "E24" "  -2.03%  "
Set-TextValue "D25" "3.652.38"
Set-TextValue "E25" "  -2.39%  "
Set-TextValue "E26" "  +0.10%  "
Set-TextValue "D27" "0.0000115"
Set-TextValue "E27" "  -5.59%  "
Set-TextValue "D28" "2.45"
Set-TextValue "E28" "  -3.00%  "
Set-TextValue "D29" "7.79"
Set-TextValue "E29" "  -4.10%  "
Set-TextValue "D30" "8.87"
Set-TextValue "E30" "  -5.21%  "
Set-TextValue "D31" "0.997"
Set-TextValue "E31" "  -0.39%  "
Set-TextValue "D32" "3.519.60"
Set-TextValue "E32" "  -2.15%  "
Set-TextValue "D33" "0.153"
Set-TextValue "E33" "  -2.10%  "
Set-TextValue "D34" "24.31"
Set-TextValue "E34" "  -4.71%  "
Set-TextValue "E35" "  -0.04%  "
Set-TextValue "B36" "Fetch.AI"
Set-TextValue "C36" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D36" "1.32"
Set-TextValue "E36" "  -9.53%  "
Set-TextValue "B37" "Aptos"
Set-TextValue "C37" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D37" "7.50"
Set-TextValue "E37" "  -4.46%  "
Set-TextValue "D38" "175.23"
Set-TextValue "E38" "  +0.34%  "
Set-TextValue "D39" "5.28"
Set-TextValue "E39" "  -6.48%  "
Set-TextValue "E40" "  -8.92%  "
Set-TextValue "D41" "0.0816"
Set-TextValue "E41" "  -4.64%  "
Set-TextValue "E42" "  -3.06%  "
Set-TextValue "D43" "0.853"
Set-TextValue "E43" "  -4.69%  "
Set-TextValue "D44" "45.27"
Set-TextValue "E44" "  -1.47%  "
Set-TextValue "E45" "  -8.42%  "
Set-TextValue "E46" "  -0.13%  "
Set-TextValue "D47" "2.36"
Set-TextValue "E47" "  -7.13%  "
Set-TextValue "D48" "23.54"
Set-TextValue "E48" "  -0.99%  "
Set-TextValue "E49" "  -2.73%  "
Set-TextValue "D50" "1.09"
Set-TextValue "E50" "  -9.37%  "
Set-TextValue "D51" "0.902"
Set-TextValue "E51" "  -4.68%  "
